$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "isbn"
$ws.Range("D1").Value = "aisle"
$ws.Range("E1").Value = "author"

# Row 2 - book1
$ws.Range("A2").Value = "book1"
$ws.Range("B2").Value = "name1"
$ws.Range("C2").Value = "isbn1"
$ws.Range("D2").Value = "aisle1"
$ws.Range("E2").Value = "author1"

# Row 3 - book2
$ws.Range("A3").Value = "book2"
$ws.Range("B3").Value = "name2"
$ws.Range("C3").Value = "isbn2"
$ws.Range("D3").Value = "aisle2"
$ws.Range("E3").Value = "author2"

# Row 4 - book3
$ws.Range("A4").Value = "book3"
$ws.Range("B4").Value = "name3"
$ws.Range("C4").Value = "isbn3"
$ws.Range("D4").Value = "aisle3"
$ws.Range("E4").Value = "author3"
